$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.390.70"
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").Value = "2.428.85"
$ws.Range("E3").Value = "  +3.61%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.98%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("D9").Value = "2.429.26"
$ws.Range("E9").Value = "  +3.72%  "
$ws.Range("E10").Value = "  +5.52%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("E13").Value = "  +4.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.86%  "
$ws.Range("E15").Value = "  +10.40%  "
$ws.Range("D16").Value = "2.869.63"
$ws.Range("D17").Value = "62.127.71"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("D18").Value = "2.429.61"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "325.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("E21").Value = "  +2.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.21%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "575.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +16.31%  "
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0948"
$ws.Range("E30").Value = "  +11.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.09%  "
$ws.Range("E33").Value = "  +2.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.19%  "
$ws.Range("E35").Value = "  +5.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.93%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  +5.94%  "
$ws.Range("E40").Value = "  +2.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "147.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.43%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("E45").Value = "  +14.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.21%  "
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0546"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.592"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("E51").Value = "  +4.78%  "
